# SEM Update on August 19th, 2025
# Adds 06/01/2025 and 07/01/2025 columns to TABLE_1 (levels) and TABLE_2 (pct change),
# and revises a handful of trailing-month values that were provisional.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TABLE_1")

# New header dates - force text so Excel does not coerce them into date serials
$ws1.Range("FG4:FH4").NumberFormat = "@"
$ws1.Range("FG4").Value = "06/01/2025"
$ws1.Range("FH4").Value = "07/01/2025"
$ws1.Range("FG4:FH4").ClearFormats()

# Revised prior-month values
$ws1.Range("FF5").Value = 11093.6
$ws1.Range("FF6").Value = 180.5
$ws1.Range("FF8").Value = 192.5
$ws1.Range("FF11").Value = 243.8
$ws1.Range("FF15").Value = 489.9
$ws1.Range("FF16").Value = 354.9
$ws1.Range("FF18").Value = 61.5
$ws1.Range("FF20").Value = 234.3
$ws1.Range("FF21").Value = 149.8
$ws1.Range("FF22").Value = 134.7
$ws1.Range("FF23").Value = 165.9
$ws1.Range("FF29").Value = 210.4
$ws1.Range("FF33").Value = 97.8
$ws1.Range("FF34").Value = 90.8
$ws1.Range("FF35").Value = 44.4
$ws1.Range("FF36").Value = 332.7
$ws1.Range("FF41").Value = 382.4
$ws1.Range("FF46").Value = 166.9
$ws1.Range("FF48").Value = 214.3
$ws1.Range("FF49").Value = 1223.5
$ws1.Range("FF53").Value = 266.2
$ws1.Range("FF55").Value = 205.5

# New 06/01/2025 and 07/01/2025 values
$ws1.Range("FG5").Value = 10506.5
$ws1.Range("FH5").Value = 9339.4
$ws1.Range("FG6").Value = 177.3
$ws1.Range("FH6").Value = 167.5
$ws1.Range("FG7").Value = 22.8
$ws1.Range("FH7").Value = 18.8
$ws1.Range("FG8").Value = 150.8
$ws1.Range("FH8").Value = 147.3
$ws1.Range("FG9").Value = 94.5
$ws1.Range("FH9").Value = 83.4
$ws1.Range("FG10").Value = 1315.5
$ws1.Range("FH10").Value = 1151.9
$ws1.Range("FG11").Value = 228.4
$ws1.Range("FH11").Value = 214.6
$ws1.Range("FG12").Value = 117.9
$ws1.Range("FH12").Value = 99.7
$ws1.Range("FG13").Value = 35.7
$ws1.Range("FH13").Value = 33.9
$ws1.Range("FG15").Value = 409.2
$ws1.Range("FH15").Value = 400.7
$ws1.Range("FG16").Value = 342.2
$ws1.Range("FH16").Value = 328.1
$ws1.Range("FG18").Value = 59.2
$ws1.Range("FH18").Value = 53.2
$ws1.Range("FG19").Value = 430
$ws1.Range("FH19").Value = 396.8
$ws1.Range("FG20").Value = 196.6
$ws1.Range("FH20").Value = 181
$ws1.Range("FG21").Value = 137.1
$ws1.Range("FH21").Value = 117
$ws1.Range("FG22").Value = 120
$ws1.Range("FH22").Value = 98.2
$ws1.Range("FG23").Value = 155.9
$ws1.Range("FH23").Value = 131.2
$ws1.Range("FG24").Value = 140.5
$ws1.Range("FH24").Value = 133.3
$ws1.Range("FG25").Value = 47.2
$ws1.Range("FH25").Value = 39.7
$ws1.Range("FG26").Value = 225.5
$ws1.Range("FH26").Value = 213.2
$ws1.Range("FG27").Value = 241
$ws1.Range("FH27").Value = 213.9
$ws1.Range("FG28").Value = 293.4
$ws1.Range("FH28").Value = 265.8
$ws1.Range("FG29").Value = 202.3
$ws1.Range("FH29").Value = 168.9
$ws1.Range("FG30").Value = 99.2
$ws1.Range("FH30").Value = 93.7
$ws1.Range("FG32").Value = 41.8
$ws1.Range("FH32").Value = 35
$ws1.Range("FG33").Value = 90.8
$ws1.Range("FH33").Value = 80.5
$ws1.Range("FG34").Value = 75.9
$ws1.Range("FH34").Value = 73.1
$ws1.Range("FG35").Value = 40.3
$ws1.Range("FH35").Value = 33.1
$ws1.Range("FG36").Value = 325.8
$ws1.Range("FH36").Value = 263.6
$ws1.Range("FG37").Value = 73.9
$ws1.Range("FH37").Value = 65.6
$ws1.Range("FG38").Value = 662.3
$ws1.Range("FH38").Value = 555.1
$ws1.Range("FG39").Value = 324.6
$ws1.Range("FH39").Value = 267.1
$ws1.Range("FG40").Value = 35.2
$ws1.Range("FH40").Value = 29
$ws1.Range("FG41").Value = 358
$ws1.Range("FH41").Value = 349.5
$ws1.Range("FG42").Value = 142.4
$ws1.Range("FH42").Value = 129.2
$ws1.Range("FG43").Value = 143.6
$ws1.Range("FH43").Value = 112.2
$ws1.Range("FG44").Value = 307.1
$ws1.Range("FH44").Value = 279.8
$ws1.Range("FG45").Value = 26.6
$ws1.Range("FH45").Value = 23.1
$ws1.Range("FG46").Value = 161.4
$ws1.Range("FH46").Value = 151.7
$ws1.Range("FG47").Value = 34.7
$ws1.Range("FH47").Value = 28.5
$ws1.Range("FG48").Value = 201.1
$ws1.Range("FH48").Value = 176.2
$ws1.Range("FG49").Value = 1175.5
$ws1.Range("FH49").Value = 1094.2
$ws1.Range("FG50").Value = 130.7
$ws1.Range("FH50").Value = 114.8
$ws1.Range("FG51").Value = 30
$ws1.Range("FH51").Value = 25.7
$ws1.Range("FG52").Value = 322.2
$ws1.Range("FH52").Value = 294
$ws1.Range("FG53").Value = 261.2
$ws1.Range("FH53").Value = 239.5
$ws1.Range("FG54").Value = 57
$ws1.Range("FH54").Value = 52.1
$ws1.Range("FG55").Value = 187.8
$ws1.Range("FH55").Value = 173.1
$ws1.Range("FG56").Value = 29
$ws1.Range("FH56").Value = 25

$ws2 = $wb.Worksheets.Item("TABLE_2")

# New header dates - force text so Excel does not coerce them into date serials
$ws2.Range("EU4:EV4").NumberFormat = "@"
$ws2.Range("EU4").Value = "06/01/2025"
$ws2.Range("EV4").Value = "07/01/2025"
$ws2.Range("EU4:EV4").ClearFormats()

# Revised prior-month values
$ws2.Range("ET5").Value = 1.00241270997406
$ws2.Range("ET6").Value = 1.06382978723403
$ws2.Range("ET8").Value = -1.88583078491335
$ws2.Range("ET11").Value = 4.32178005990587
$ws2.Range("ET15").Value = 0.0408413314274138
$ws2.Range("ET16").Value = 1.28424657534248
$ws2.Range("ET18").Value = -0.485436893203879
$ws2.Range("ET20").Value = -1.47182506307822
$ws2.Range("ET21").Value = 1.1478730587441
$ws2.Range("ET22").Value = 1.81405895691608
$ws2.Range("ET23").Value = 1.28205128205126
$ws2.Range("ET29").Value = -0.142382534409118
$ws2.Range("ET33").Value = 1.66320166320167
$ws2.Range("ET34").Value = 4.60829493087558
$ws2.Range("ET35").Value = 3.25581395348839
$ws2.Range("ET36").Value = 1.64986251145739
$ws2.Range("ET41").Value = -1.7724120215772
$ws2.Range("ET46").Value = 0.481637567730273
$ws2.Range("ET48").Value = 3.32690453230474
$ws2.Range("ET49").Value = 1.25796573698585
$ws2.Range("ET53").Value = -1.51683314835367
$ws2.Range("ET55").Value = 0.488997555012225

# New 06/01/2025 and 07/01/2025 values
$ws2.Range("EU5").Value = 1.59453082695134
$ws2.Range("EV5").Value = 0.733438315680499
$ws2.Range("EU6").Value = 0.853242320819113
$ws2.Range("EV6").Value = 1.08630054315028
$ws2.Range("EU7").Value = 2.70270270270271
$ws2.Range("EV7").Value = -3.09278350515463
$ws2.Range("EU8").Value = -3.33333333333333
$ws2.Range("EV8").Value = -2.77227722772276
$ws2.Range("EU9").Value = 0.638977635782741
$ws2.Range("EV9").Value = 0.60313630880579
$ws2.Range("EU10").Value = 2.82163514147256
$ws2.Range("EV10").Value = 3.44858554108667
$ws2.Range("EU11").Value = 4.53089244851259
$ws2.Range("EV11").Value = 4.78515624999999
$ws2.Range("EU12").Value = 0.16992353440953
$ws2.Range("EV12").Value = 0.402819738167147
$ws2.Range("EU13").Value = -1.65289256198346
$ws2.Range("EV13").Value = -0.294117647058828
$ws2.Range("EU15").Value = 0.0488997555012336
$ws2.Range("EV15").Value = 0.124937531234369
$ws2.Range("EU16").Value = 2.39377618192697
$ws2.Range("EV16").Value = 2.65957446808511
$ws2.Range("EU18").Value = -0.504201680672264
$ws2.Range("EV18").Value = 0.188323917137466
$ws2.Range("EU19").Value = 2.8462090408993
$ws2.Range("EV19").Value = -0.800000000000011
$ws2.Range("EU20").Value = -1.15635997988938
$ws2.Range("EV20").Value = 5.72429906542057
$ws2.Range("EU21").Value = 3.54984894259822
$ws2.Range("EV21").Value = 1.38648180242634
$ws2.Range("EU22").Value = 1.09519797809605
$ws2.Range("EV22").Value = 0.61475409836065
$ws2.Range("EU23").Value = 1.16807268007788
$ws2.Range("EV23").Value = 1.46945088940447
$ws2.Range("EU24").Value = 0.861450107681255
$ws2.Range("EV24").Value = 0.755857898715041
$ws2.Range("EU25").Value = 1.50537634408603
$ws2.Range("EV25").Value = 0.761421319796965
$ws2.Range("EU26").Value = 3.72585096596137
$ws2.Range("EV26").Value = 3.79746835443039
$ws2.Range("EU27").Value = -0.248344370860937
$ws2.Range("EV27").Value = -0.418994413407824
$ws2.Range("EU28").Value = 0.548320767649063
$ws2.Range("EV28").Value = 1.5666794038976
$ws2.Range("EU29").Value = 1.14999999999999
$ws2.Range("EV29").Value = -1.57342657342657
$ws2.Range("EU30").Value = 1.43149284253576
$ws2.Range("EV30").Value = -1.57563025210086
$ws2.Range("EU32").Value = -1.41509433962264
$ws2.Range("EV32").Value = -1.9607843137255
$ws2.Range("EU33").Value = 2.83125707814269
$ws2.Range("EV33").Value = 1.64141414141416
$ws2.Range("EU34").Value = 3.54706684856752
$ws2.Range("EV34").Value = 3.68794326241134
$ws2.Range("EU35").Value = 2.8061224489796
$ws2.Range("EV35").Value = 3.11526479750779
$ws2.Range("EU36").Value = 1.84432635198501
$ws2.Range("EV36").Value = -0.715630885122402
$ws2.Range("EU37").Value = 3.35664335664336
$ws2.Range("EV37").Value = 3.63349131121643
$ws2.Range("EU38").Value = 0.760687661646145
$ws2.Range("EV38").Value = 1.01910828025478
$ws2.Range("EU39").Value = 1.43750000000001
$ws2.Range("EV39").Value = 1.71363290175171
$ws2.Range("EU40").Value = 2.32558139534885
$ws2.Range("EV40").Value = 0.34602076124568
$ws2.Range("EU41").Value = -1.18686171680928
$ws2.Range("EV41").Value = 0.402183280666469
$ws2.Range("EU42").Value = 2.89017341040465
$ws2.Range("EV42").Value = 2.05371248025276
$ws2.Range("EU43").Value = 1.41242937853109
$ws2.Range("EV43").Value = -1.49253731343284
$ws2.Range("EU44").Value = 1.85737976782754
$ws2.Range("EV44").Value = 1.56079854809438
$ws2.Range("EU45").Value = -5.33807829181494
$ws2.Range("EV45").Value = -2.53164556962026
$ws2.Range("EU46").Value = 0.560747663551405
$ws2.Range("EV46").Value = 1.47157190635451
$ws2.Range("EU47").Value = 1.16618075801749
$ws2.Range("EV47").Value = -0.696864111498255
$ws2.Range("EU48").Value = 3.07534597642233
$ws2.Range("EV48").Value = 3.52526439482961
$ws2.Range("EU49").Value = 1.24020325553355
$ws2.Range("EV49").Value = 0.82004975582789
$ws2.Range("EU50").Value = 3.64789849325931
$ws2.Range("EV50").Value = 4.26884650317894
$ws2.Range("EU51").Value = 0
$ws2.Range("EV51").Value = -0.387596899224798
$ws2.Range("EU52").Value = 1.8975332068311
$ws2.Range("EV52").Value = 4.62633451957295
$ws2.Range("EU53").Value = -0.985595147839259
$ws2.Range("EV53").Value = -2.20498162515313
$ws2.Range("EU54").Value = 1.24333925399645
$ws2.Range("EV54").Value = 1.16504854368931
$ws2.Range("EU55").Value = 0.374131480491725
$ws2.Range("EV55").Value = 1.10981308411215
$ws2.Range("EU56").Value = 1.39860139860139
$ws2.Range("EV56").Value = 3.73443983402491
